# Re-run Solver on the "Optimization Model" sheet producing a second set of
# reports (Answer Report 2, Sensitivity Report 2, Limits Report 2), inserted
# before the original "Optimization Model" sheet, and update the model
# itself with the new solution values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the three report sheets, placing the copies in front of
#    "Optimization Model", and rename them to the "2" variants.
# ---------------------------------------------------------------------
$target = $wb.Worksheets.Item("Optimization Model")
$wb.Worksheets.Item("Answer Report 1").Copy($target)
$wb.Worksheets.Item(1).Name = "Answer Report 2"

$target = $wb.Worksheets.Item("Optimization Model")
$wb.Worksheets.Item("Sensitivity Report 1").Copy($target)
$wb.Worksheets.Item(2).Name = "Sensitivity Report 2"

$target = $wb.Worksheets.Item("Optimization Model")
$wb.Worksheets.Item("Limits Report 1").Copy($target)
$wb.Worksheets.Item(3).Name = "Limits Report 2"

# ---------------------------------------------------------------------
# 2. Update "Answer Report 2" with the new solver run's numbers.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Answer Report 2")
$ws.Range("A2").Value = "Worksheet: [Excel-Solver-Example-LP.xlsx]Optimization Model"
$ws.Range("A3").Value = "Report Created: 1/27/22 10:37:00 PM"
$ws.Range("B7").Value = "Solution Time: 4295291.228 Seconds."
$ws.Range("B8").Value = "Iterations: 3 Subproblems: 0"
$ws.Range("B10").Value = "Max Time Unlimited, Iterations Unlimited, Precision 0.000001"

$ws.Range("D16").Value = 172706.33333333331
$ws.Range("E16").Value = 252000

$ws.Range("D22").Value = 11336.666666666666
$ws.Range("E22").Value = 18000
$ws.Range("D23").Value = 666.33333333333326
$ws.Range("E23").Value = 0
$ws.Range("D24").Value = 16004
$ws.Range("E24").Value = 24000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 4001
$ws.Range("E26").Value = 6000

$ws.Range("F32").Value = "Not Binding"
$ws.Range("G32").Value = 1

# ---------------------------------------------------------------------
# 3. Update "Sensitivity Report 2".
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sensitivity Report 2")
$ws.Range("A2").Value = "Worksheet: [Excel-Solver-Example-LP.xlsx]Optimization Model"
$ws.Range("A3").Value = "Report Created: 1/27/22 10:37:01 PM"

$ws.Range("D10").Value = 18000
$ws.Range("G10").Value = [double]"1E+30"
$ws.Range("H10").Value = 11.9
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = -119
$ws.Range("H11").Value = [double]"1E+30"
$ws.Range("D12").Value = 24000
$ws.Range("H12").Value = 8.8888888888888893
$ws.Range("E13").Value = -80
$ws.Range("G13").Value = 80
$ws.Range("D14").Value = 6000

$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4001
$ws.Range("G20").Value = [double]"1E+30"
$ws.Range("H20").Value = 1
$ws.Range("E21").Value = 42
$ws.Range("G21").Value = [double]"1E+30"
$ws.Range("H21").Value = 6000
$ws.Range("H22").Value = 18000
$ws.Range("H23").Value = 24000

# ---------------------------------------------------------------------
# 4. Update "Limits Report 2".
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Limits Report 2")
$ws.Range("A2").Value = "Worksheet: [Excel-Solver-Example-LP.xlsx]Optimization Model"
$ws.Range("A3").Value = "Report Created: 1/27/22 10:37:01 PM"

$ws.Range("D8").Value = 252000
$ws.Range("D14").Value = 18000
$ws.Range("D15").Value = 0
$ws.Range("D16").Value = 24000
$ws.Range("D18").Value = 6000

# ---------------------------------------------------------------------
# 5. Update the "Optimization Model" sheet itself with the new solution.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Optimization Model")
$ws.Range("C8").Value = 18000
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 24000
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 6000

$ws.Range("H11").Value = 4000
$ws.Range("J11").Value = 4001

# ---------------------------------------------------------------------
# 6. Fix up selections / active sheet to match the new state.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Sensitivity Report 2").Range("A1:A3").Select()
$wb.Worksheets.Item("Limits Report 2").Range("A1:A3").Select()

$wb.Worksheets.Item("Sensitivity Report 1").Range("E20").Select()

$ws = $wb.Worksheets.Item("Optimization Model")
$ws.Range("B15").Select()
$ws.Activate()
